$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 551; existing rows 551..571 shift down to 552..572
$ws.Rows.Item(551).Insert()

# Populate the newly inserted row 551 with the new weekly observation
$ws.Cells.Item(551, 1).Value = 3
$ws.Cells.Item(551, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(551, 3).Value = "Coquimbo"
$ws.Cells.Item(551, 4).Value = 45075
$ws.Cells.Item(551, 5).Value = 5
$ws.Cells.Item(551, 6).Value = 100112040
$ws.Cells.Item(551, 7).Value = "Cilantro"
$ws.Cells.Item(551, 8).Value = "Sin especificar"
$ws.Cells.Item(551, 9).Value = "Primera"
$ws.Cells.Item(551, 10).Value = 120
$ws.Cells.Item(551, 11).Value = 4500
$ws.Cells.Item(551, 12).Value = 4500
$ws.Cells.Item(551, 13).Value = 4500
$ws.Cells.Item(551, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(551, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(551, 16).Value = 1500
$ws.Cells.Item(551, 17).Value = 3
$ws.Cells.Item(551, 18).Value = "Hortaliza"
